$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections across rows 1-8 (pre/post/total fixation metrics re-split) ---
# Only the cells whose values actually change are touched, so every other
# cell (including the already-empty placeholder cells) is left exactly as-is.

# Row 1
$ws.Range("A1").Value = ""
# Row 2
$ws.Range("A2").Value = 'Fixation'
$ws.Range("B2").Value = 'based'
$ws.Range("C2").Value = 'metrics'
# Row 3
$ws.Range("A3").Value = 'Revisit'
$ws.Range("B3").Value = 31
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 16
$ws.Range("H3").Value = 21
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 2
$ws.Range("Q3").Value = 0
$ws.Range("S3").Value = 12
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 8
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 29
$ws.Range("AD3").Value = ""
$ws.Range("AH3").Value = 0
$ws.Range("AJ3").Value = 3
$ws.Range("AO3").Value = 4
$ws.Range("AR3").Value = ""
# Row 4
$ws.Range("A4").Value = 'Fixation'
$ws.Range("B4").Value = 112
$ws.Range("C4").Value = 83
$ws.Range("D4").Value = 17
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 26
$ws.Range("H4").Value = 38
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 5
$ws.Range("Q4").Value = 1
$ws.Range("S4").Value = 15
$ws.Range("T4").Value = 4
$ws.Range("U4").Value = 4
$ws.Range("V4").Value = 13
$ws.Range("W4").Value = 9
$ws.Range("X4").Value = 123
$ws.Range("AD4").Value = ""
$ws.Range("AH4").Value = 2
$ws.Range("AJ4").Value = 5
$ws.Range("AO4").Value = 5
$ws.Range("AR4").Value = ""
# Row 5
$ws.Range("A5").Value = 'Dwell'
$ws.Range("B5").Value = 48876
$ws.Range("C5").Value = 40951.13
$ws.Range("D5").Value = 7240.69
$ws.Range("E5").Value = 2936.44
$ws.Range("F5").Value = 3637.27
$ws.Range("G5").Value = 15825.51
$ws.Range("H5").Value = 22098.39
$ws.Range("J5").Value = 1017.74
$ws.Range("K5").Value = 1985.38
$ws.Range("L5").Value = 3236.64
$ws.Range("Q5").Value = 408.83
$ws.Range("S5").Value = 6848.81
$ws.Range("T5").Value = 934.35
$ws.Range("U5").Value = 1251.27
$ws.Range("V5").Value = 8917.97
$ws.Range("W5").Value = 5163.77
$ws.Range("X5").Value = 75329.62
$ws.Range("AD5").Value = ""
$ws.Range("AH5").Value = 417.12
$ws.Range("AJ5").Value = 1610.02
$ws.Range("AO5").Value = 1718.51
$ws.Range("AR5").Value = ""
# Row 6
$ws.Range("A6").Value = 'Dwell'
$ws.Range("B6").Value = 27.67
$ws.Range("C6").Value = 23.18
$ws.Range("D6").Value = 4.1
$ws.Range("E6").Value = 1.66
$ws.Range("F6").Value = 2.06
$ws.Range("G6").Value = 8.96
$ws.Range("H6").Value = 12.51
$ws.Range("J6").Value = 0.58
$ws.Range("K6").Value = 1.12
$ws.Range("L6").Value = 1.83
$ws.Range("N6").Value = 1.68
$ws.Range("P6").Value = 0.68
$ws.Range("Q6").Value = 0.23
$ws.Range("S6").Value = 3.88
$ws.Range("T6").Value = 0.53
$ws.Range("U6").Value = 0.71
$ws.Range("V6").Value = 5.05
$ws.Range("W6").Value = 2.92
$ws.Range("X6").Value = 42.64
$ws.Range("Y6").Value = 4.59
$ws.Range("Z6").Value = 1.12
$ws.Range("AA6").Value = 0.6
$ws.Range("AB6").Value = 1.12
$ws.Range("AC6").Value = 0.6
$ws.Range("AD6").Value = ""
$ws.Range("AG6").Value = 0.34
$ws.Range("AH6").Value = 0.24
$ws.Range("AJ6").Value = 0.91
$ws.Range("AL6").Value = 5.53
$ws.Range("AM6").Value = 0.59
$ws.Range("AN6").Value = 0.58
$ws.Range("AO6").Value = 0.97
$ws.Range("AR6").Value = ""
# Row 7
$ws.Range("A7").Value = 'Fixation'
$ws.Range("B7").Value = 436.39
$ws.Range("C7").Value = 493.39
$ws.Range("D7").Value = 425.92
$ws.Range("E7").Value = 489.41
$ws.Range("F7").Value = 909.32
$ws.Range("G7").Value = 608.67
$ws.Range("H7").Value = 581.54
$ws.Range("J7").Value = 339.25
$ws.Range("K7").Value = 1985.38
$ws.Range("L7").Value = 647.33
$ws.Range("Q7").Value = 408.83
$ws.Range("S7").Value = 456.59
$ws.Range("T7").Value = 233.59
$ws.Range("U7").Value = 312.82
$ws.Range("V7").Value = 686
$ws.Range("W7").Value = 573.75
$ws.Range("X7").Value = 612.44
$ws.Range("AD7").Value = ""
$ws.Range("AH7").Value = 208.56
$ws.Range("AJ7").Value = 322
$ws.Range("AO7").Value = 343.7
$ws.Range("AR7").Value = ""
# Row 8
$ws.Range("A8").Value = 'First'
$ws.Range("B8").Value = 350.25
$ws.Range("C8").Value = 350.25
$ws.Range("L8").Value = 350.25
$ws.Range("V8").Value = 350.25
$ws.Range("AD8").Value = ""
$ws.Range("AR8").Value = ""

# --- Structural trim: drop the three trailing blank rows (9-11) ---
$ws.Rows("9:11").Delete()

# --- Header row formatting: the bold / thin-border / centered style (cellXfs idx 1) is gone ---
$ws.Range("A1:AR1").ClearFormats()
